# This script applies the "generates random code and adds to database" edit:
#  - appends newly-purchased games to a few users' Library column on accountInfo
#  - adds a new "Valid_Codes" column (G) to the store sheet with generated codes
#    for a handful of titles
#  - updates the active sheet/selection to reflect the editor's final cursor position

$wb = $excel.ActiveWorkbook
$wsAccount = $wb.Worksheets.Item("accountInfo")
$wsStore   = $wb.Worksheets.Item("store")

# ---------------------------------------------------------------------------
# accountInfo sheet: extend each user's Library (column D) with new purchases
# ---------------------------------------------------------------------------

# Bob (row 2) buys several more games
$bobLibrary = $wsAccount.Range("D2").Value2
$wsAccount.Range("D2").Value = $bobLibrary + "/Stardew Valley/Portal/RollerCoaster Tycoon/Team Fortress 2/The Sims/RollerCoaster Tycoon 2"

# Amy (row 3) buys a couple more games
$amyLibrary = $wsAccount.Range("D3").Value2
$wsAccount.Range("D3").Value = $amyLibrary + "/The Witcher 3: Wild Hunt/Celeste"

# Rachel (row 5) didn't have a Library cell before - she now has one
$wsAccount.Range("D5").Value = "Civilization VI/The Witcher 3: Wild Hunt/Celeste/Doom/Stardew Valley/Cuphead/Portal"

# Jose (row 7) didn't have a Library cell before - he now has one
$wsAccount.Range("D7").Value = "Age of Empires/Doom/The Witcher 3: Wild Hunt/Portal/Celeste/Half-Life: Alyx/Stardew Valley/Myst/Tetris Effect/RollerCoaster Tycoon/The Sims/Minecraft"

# John (row 12) had an empty, numerically-formatted Library cell - give him a library
$wsAccount.Range("D12").Value = "Fallout/The Witcher 3: Wild Hunt/The Sims"

# ---------------------------------------------------------------------------
# store sheet: add a new "Valid_Codes" column (G) with randomly generated codes
# ---------------------------------------------------------------------------

# Header cell G1 - copy formatting from F1 (the other bold/header-styled column)
# so the new column header matches the existing header style, then set its text.
$wsStore.Range("F1").Copy()
$wsStore.Range("G1").PasteSpecial(-4122)
$wsStore.Range("G1").Value = "Valid_Codes"

# Randomly generated redemption codes for a handful of titles
$wsStore.Range("G3").Value = "087QB0A/OJBNA1W"
$wsStore.Range("G4").Value = "D2WKY6W"
$wsStore.Range("G5").Value = "AVDF58I/PFCG169"
$wsStore.Range("G7").Value = "Y6BAP56"

# ---------------------------------------------------------------------------
# Final cursor/selection state: store sheet becomes active, cell G2 selected
# ---------------------------------------------------------------------------
$wsStore.Activate() | Out-Null
$wsStore.Range("G2").Select() | Out-Null
